# modificacion en metodo upload database
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")

# Fix typo'd marca value "lllll" -> "asd"
$ws.Range("B2").Value = "asd"

# distribucion_tiendas column: set all data rows (2-24) to 1
for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 4).Value = 1
}
